$wb = $excel.ActiveWorkbook

$newVersion = "Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on February 03 2026 17.29.55 EST)"

# --- Update the "About" sheet ---
$aboutWs = $wb.Worksheets.Item("About")

$aboutWs.Range("A2").Value = "Version: " + $newVersion

$newCitation = "Recommended Citation:  `"Global Energy Monitor, Coal mine boundaries and methane sources for Tahmoor Coal Mine, Australia, M0103, version '" + $newVersion + "'. (See the CC license for attribution requirements if sharing or adapting the data set.)"

$aboutWs.Range("A6").Value = $newCitation

# --- Update the "Boundaries and methane sources" sheet ---
$dataWs = $wb.Worksheets.Item("Boundaries and methane sources")

for ($row = 2; $row -le 8; $row++) {
    $dataWs.Cells.Item($row, 19).Value = $newVersion  # column S = 19 = build_version
}
